$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct previously-reported AgTests (column H) values for several existing rows ---
$ws.Range("H223").Value = 1152
$ws.Range("H230").Value = 614
$ws.Range("H232").Value = 2138
$ws.Range("H239").Value = 6185
$ws.Range("H240").Value = 41132
$ws.Range("H245").Value = 3596
$ws.Range("H246").Value = 2114
$ws.Range("H248").Value = 45188
$ws.Range("H253").Value = 5636
$ws.Range("H254").Value = 6145
$ws.Range("H256").Value = 1245
$ws.Range("H258").Value = 3859
$ws.Range("H259").Value = 6104
$ws.Range("H260").Value = 11611
$ws.Range("H267").Value = 13242
$ws.Range("H271").Value = 42584
$ws.Range("H275").Value = 28741
$ws.Range("H278").Value = 30001
$ws.Range("H279").Value = 43382
$ws.Range("H280").Value = 35570
$ws.Range("H281").Value = 45201
$ws.Range("H282").Value = 46781
$ws.Range("H285").Value = 40884
$ws.Range("H286").Value = 54196
$ws.Range("H287").Value = 57619
$ws.Range("H288").Value = 56086
$ws.Range("H289").Value = 64384
$ws.Range("H291").Value = 14831
$ws.Range("H292").Value = 81212
$ws.Range("H293").Value = 81666
$ws.Range("H294").Value = 90790
$ws.Range("H295").Value = 19353
$ws.Range("H300").Value = 70476
$ws.Range("H302").Value = 72555
$ws.Range("H306").Value = 70725
$ws.Range("H307").Value = 73347
$ws.Range("H309").Value = 57092
$ws.Range("H310").Value = 90810
$ws.Range("H313").Value = 72891
$ws.Range("H314").Value = 65070
$ws.Range("H316").Value = 49040
$ws.Range("H318").Value = 24467
$ws.Range("H320").Value = 83746
$ws.Range("H322").Value = 98713

# --- Append new daily record for 2021-01-23 update (date serial 44218) ---
$ws.Range("A324").Value = 44218
$ws.Range("B324").Value = 50063
$ws.Range("C324").Value = 186683
$ws.Range("D324").Value = -140585
$ws.Range("E324").Value = 20264
$ws.Range("F324").Value = -182964
$ws.Range("G324").Value = 3965
$ws.Range("H324").Value = 30746
$ws.Range("I324").Value = 437
